$d = $word.ActiveDocument

$replacements = @(
    @("500×6=3000", "374×3=1122"),
    @("393×9=3537", "104×3=312"),
    @("626×8=5008", "198×3=594"),
    @("368×4=1472", "251×5=1255"),
    @("766×2=1532", "463×5=2315"),
    @("820×7=5740", "550×5=2750"),
    @("139×2=278", "387×5=1935"),
    @("586×6=3516", "554×6=3324"),
    @("876×9=7884", "604×7=4228"),
    @("635×5=3175", "783×7=5481"),
    @("873×7=6111", "974×5=4870"),
    @("455×4=1820", "110×2=220"),
    @("504×5=2520", "311×4=1244"),
    @("150×8=1200", "115×4=460"),
    @("309×7=2163", "542×6=3252"),
    @("933×9=8397", "964×9=8676"),
    @("721×5=3605", "414×3=1242"),
    @("761×8=6088", "230×7=1610"),
    @("611×8=4888", "826×5=4130"),
    @("721×7=5047", "120×5=600"),
    @("223×5=1115", "419×8=3352"),
    @("928×7=6496", "132×3=396"),
    @("365×2=730", "554×4=2216"),
    @("448×5=2240", "837×4=3348"),
    @("354×5=1770", "555×3=1665")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
